# Update row 2 of the cash-flow data sheet with the latest reporting period.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2 holds a zero-padded numeric-looking code ("001" -> "004"). Force it to
# stay text (matching the other REPORT_TYPE_CODE-style columns) instead of
# letting it collapse to the number 4, then restore the default style so we
# don't leave a stray number format behind on the cell.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "004"
$ws.Range("J2").Style = "Normal"

# K2 is unchanged ("001") - left untouched.

# REPORT_DATE moves from 2019-12-31 to 2020-09-30.
$ws.Range("N2").Value = "2020-09-30 00:00:00"

# Updated cash-flow figures / ratios for the new reporting period.
$ws.Range("O2").Value = 70497982.62
$ws.Range("P2").Value = 267.3711558162
$ws.Range("Q2").Value = 312450800.79
$ws.Range("R2").Value = 1185.0031538241
$ws.Range("S2").Value = 39133838.92
$ws.Range("T2").Value = 148.4192788887
$ws.Range("U2").Value = -72107221.13
$ws.Range("V2").Value = -273.4743653609
$ws.Range("W2").Value = 1642791.79
$ws.Range("X2").Value = 6.2304639556
$ws.Range("Y2").Value = 73323303.28
$ws.Range("Z2").Value = 278.0864872675
$ws.Range("AA2").Value = 28277154.26
$ws.Range("AB2").Value = 107.2441385797
$ws.Range("AC2").Value = 26367086.01
$ws.Range("AD2").Value = 1518.5581256694
